# The spreadsheet used to ship with a 20-row "blank template" block
# (rows 5-24) styled as a bordered entry grid, below the 3 sample data
# rows. That fixed-size block is no longer needed, so remove it,
# restoring the sheet to just the header + the 3 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused template rows entirely (shifts nothing below them
# back up since they were the last rows on the sheet).
$ws.Rows("5:24").Delete() | Out-Null

# Row 4 had a custom row height + thick bottom border because it used to
# sit right above the bordered block; with that block gone it should go
# back to the sheet's normal auto-sized row formatting.
$ws.Rows(4).AutoFit() | Out-Null

# Move the selection off the old A5:G24 block (which no longer exists)
# to a plain single-cell selection, matching a manual click elsewhere on
# the sheet after trimming the rows.
$ws.Range("V18").Select() | Out-Null
